$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 20:22"

# Row 4
$ws.Cells.Item(4, 2).Value = 6608623
$ws.Cells.Item(4, 3).Value = 20463
$ws.Cells.Item(4, 4).Value = 3884654
$ws.Cells.Item(4, 5).Value = 2527160
$ws.Cells.Item(4, 7).Value = 482
$ws.Cells.Item(4, 8).Value = 196809

# Row 5
$ws.Cells.Item(5, 2).Value = 4653302
$ws.Cells.Item(5, 3).Value = 93577
$ws.Cells.Item(5, 4).Value = 3616208
$ws.Cells.Item(5, 5).Value = 959633
$ws.Cells.Item(5, 7).Value = 1157
$ws.Cells.Item(5, 8).Value = 77461

# Row 12
$ws.Cells.Item(12, 2).Value = 566326
$ws.Cells.Item(12, 3).Value = 12183
$ws.Cells.Item(12, 7).Value = 48
$ws.Cells.Item(12, 8).Value = 29747

# Row 16
$ws.Cells.Item(16, 1).Value = "Francia"
$ws.Cells.Item(16, 2).Value = 363350
$ws.Cells.Item(16, 3).Value = 9406
$ws.Cells.Item(16, 4).Value = 89059
$ws.Cells.Item(16, 5).Value = 243398
$ws.Cells.Item(16, 7).Value = 80
$ws.Cells.Item(16, 8).Value = 30893

# Row 17
$ws.Cells.Item(17, 1).Value = "Reino Unido"
$ws.Cells.Item(17, 2).Value = 361677
$ws.Cells.Item(17, 3).Value = 3539
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 7).Value = 6
$ws.Cells.Item(17, 8).Value = 41614

# Row 21
$ws.Cells.Item(21, 2).Value = 288126
$ws.Cells.Item(21, 3).Value = 1671
$ws.Cells.Item(21, 4).Value = 256524
$ws.Cells.Item(21, 5).Value = 24651
$ws.Cells.Item(21, 7).Value = 56
$ws.Cells.Item(21, 8).Value = 6951

# Row 43
$ws.Cells.Item(43, 2).Value = 81009
$ws.Cells.Item(43, 3).Value = 703
$ws.Cells.Item(43, 4).Value = 69703
$ws.Cells.Item(43, 5).Value = 8377
$ws.Cells.Item(43, 7).Value = 11
$ws.Cells.Item(43, 8).Value = 2929

# Row 51
$ws.Cells.Item(51, 1).Value = "Etiopia"
$ws.Cells.Item(51, 2).Value = 63367
$ws.Cells.Item(51, 3).Value = 789
$ws.Cells.Item(51, 4).Value = 24024
$ws.Cells.Item(51, 5).Value = 38357
$ws.Cells.Item(51, 7).Value = 12
$ws.Cells.Item(51, 8).Value = 986

# Row 52
$ws.Cells.Item(52, 1).Value = "Portugal"
$ws.Cells.Item(52, 2).Value = 62813
$ws.Cells.Item(52, 3).Value = 687
$ws.Cells.Item(52, 4).Value = 43644
$ws.Cells.Item(52, 5).Value = 17314
$ws.Cells.Item(52, 7).Value = 3
$ws.Cells.Item(52, 8).Value = 1855

# Row 72
$ws.Cells.Item(72, 2).Value = 30571
$ws.Cells.Item(72, 3).Value = 211
$ws.Cells.Item(72, 5).Value = 5426

# Row 100
$ws.Cells.Item(100, 1).Value = "Maldivas"
$ws.Cells.Item(100, 2).Value = 8990
$ws.Cells.Item(100, 3).Value = 90
$ws.Cells.Item(100, 4).Value = 6846
$ws.Cells.Item(100, 5).Value = 2113
$ws.Cells.Item(100, 8).Value = 31

# Row 101
$ws.Cells.Item(101, 1).Value = "Tayikistan"
$ws.Cells.Item(101, 2).Value = 8977
$ws.Cells.Item(101, 3).Value = 38
$ws.Cells.Item(101, 4).Value = 7747
$ws.Cells.Item(101, 5).Value = 1158
$ws.Cells.Item(101, 8).Value = 72

# Row 116
$ws.Cells.Item(116, 1).Value = "Mozambique"
$ws.Cells.Item(116, 2).Value = 4918
$ws.Cells.Item(116, 3).Value = 86
$ws.Cells.Item(116, 4).Value = 2899
$ws.Cells.Item(116, 5).Value = 1988
$ws.Cells.Item(116, 8).Value = 31

# Row 117
$ws.Cells.Item(117, 1).Value = "Congo"
$ws.Cells.Item(117, 2).Value = 4891
$ws.Cells.Item(117, 4).Value = 3887
$ws.Cells.Item(117, 5).Value = 921
$ws.Cells.Item(117, 8).Value = 83

# Row 137
$ws.Cells.Item(137, 1).Value = "Bahamas"
$ws.Cells.Item(137, 2).Value = 2814
$ws.Cells.Item(137, 3).Value = 93
$ws.Cells.Item(137, 4).Value = 1220
$ws.Cells.Item(137, 5).Value = 1529
$ws.Cells.Item(137, 8).Value = 65

# Row 138
$ws.Cells.Item(138, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(138, 2).Value = 2777
$ws.Cells.Item(138, 3).Value = 79
$ws.Cells.Item(138, 4).Value = 762
$ws.Cells.Item(138, 5).Value = 1970
$ws.Cells.Item(138, 7).Value = 2
$ws.Cells.Item(138, 8).Value = 45

# Row 139
$ws.Cells.Item(139, 1).Value = "Jordania"
$ws.Cells.Item(139, 2).Value = 2739
$ws.Cells.Item(139, 4).Value = 1981
$ws.Cells.Item(139, 5).Value = 738
$ws.Cells.Item(139, 8).Value = 20

# Row 151
$ws.Cells.Item(151, 2).Value = 2007
$ws.Cells.Item(151, 3).Value = 4
$ws.Cells.Item(151, 5).Value = 214
$ws.Cells.Item(151, 7).Value = 2
$ws.Cells.Item(151, 8).Value = 582

